$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.689.27"
$ws.Range("E2").Value = "  +1.09%  "
$ws.Range("D3").Value = "3.799.55"
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.31%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +1.19%  "
$ws.Range("E9").Value = "  +2.41%  "
$ws.Range("E10").Value = "  -0.90%  "
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.96%  "
$ws.Range("D14").Value = "4.440.42"
$ws.Range("E14").Value = "  +1.42%  "
$ws.Range("D15").Value = "3.805.33"
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.56"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.31%  "
$ws.Range("D17").Value = "67.667.53"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("E18").Value = "  +2.88%  "
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("E20").Value = "  +1.71%  "
$ws.Range("E21").Value = "  -4.90%  "
$ws.Range("E22").Value = "  +1.48%  "
$ws.Range("E23").Value = "  +4.75%  "
$ws.Range("E24").Value = "  +1.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.95%  "
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.90%  "
$ws.Range("D29").Value = "3.939.97"
$ws.Range("E29").Value = "  +1.14%  "
$ws.Range("E30").Value = "  +0.57%  "
$ws.Range("E31").Value = "  +3.27%  "
$ws.Range("E32").Value = "  +1.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.56"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("E36").Value = "  +1.11%  "
$ws.Range("E37").Value = "  +3.98%  "
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("E39").Value = "  +1.40%  "
$ws.Range("E40").Value = "  +1.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.25%  "
$ws.Range("E44").Value = "  +1.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "147.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "394.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.18%  "
$ws.Range("E49").Value = "  +1.37%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.02%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "26.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.32%  "
